# Applies the "testee" -> "thecomp" rebrand, updated numeric stats, new
# logo_description column on the Notes sheet, and refreshed Bio / Description
# text, as described in the commit "a stable version trying to fix excel".

$wb = $excel.ActiveWorkbook

$companies = $wb.Worksheets.Item("Companies")
$bios      = $wb.Worksheets.Item("Bios")
$notes     = $wb.Worksheets.Item("Notes")

# ----------------------------------------------------------------------
# Companies sheet
# ----------------------------------------------------------------------
$companies.Range("A2").Value = 104
$companies.Range("B2").Value = "thecomp"
$companies.Range("D2").Value = "www.thecomp.com"
$companies.Range("I2").Value = "thecomp.jpg"
$companies.Range("J2").Value = "thecompBD.jpg"
$companies.Range("K2").Value = "thecompBanner.jpg"
$companies.Range("M2").Value = 19
$companies.Range("R2").Value = 41

# ----------------------------------------------------------------------
# Bios sheet
# ----------------------------------------------------------------------
$newBio = @'
Name: thecomp
Promotion Type: Professional Wrestling Company
Size: Medium
Location: TBD
Description:
thecomp is a professional wrestling company that prides itself on providing high-quality, entertaining wrestling matches for fans all over the world. As a test company, thecomp is constantly striving to push the boundaries of traditional wrestling and create new and innovative content for its audience.
Founded by a group of wrestling enthusiasts, thecomp has quickly gained a reputation for showcasing some of the best talent in the industry. With a roster of skilled and passionate wrestlers, thecomp offers a diverse range of styles and personalities, ensuring that there is something for every fan to enjoy.
In addition to its in-ring action, thecomp also focuses on creating engaging storylines and rivalries that captivate audiences and keep them coming back for more. The company also places a strong emphasis on fan interaction, regularly hosting meet-and-greet events and Q&A sessions with its wrestlers.
With a dedicated and hardworking team behind the scenes, thecomp is committed to providing a professional and entertaining product that fans can be proud to support. Whether you're a longtime wrestling fan or new to the sport, thecomp has something for everyone to enjoy. Stay tuned for upcoming events and shows from thecomp!
'@

$bios.Range("A2").Value = 104
$bios.Range("B2").Value = $newBio
# The long, multi-line bio text makes the engine auto-expand the row height;
# re-run AutoFit so the row keeps its original (default) height metadata.
$bios.Range("B2").EntireRow.AutoFit()

# ----------------------------------------------------------------------
# Notes sheet
# ----------------------------------------------------------------------
$notes.Range("A2").Value = "thecomp"
$notes.Range("B2").Value = "a test company"
$notes.Range("D2").Value = "thecomp.jpg"
$notes.Range("E2").Value = "thecompBD.jpg"
$notes.Range("F2").Value = "thecompBanner.jpg"

# New "logo_description" column (H), header styled like the rest of row 1.
$notes.Range("G1").Copy()
$notes.Range("H1").PasteSpecial(-4122)
$notes.Range("H1").Value = "logo_description"
$notes.Range("H2").Value = "The logo for 'thecomp' may feature a bold, modern font in black and white with a stylized wrestling ring silhouette in the background, symbolizing strength and competition."

Write-Host "Edit complete"
